$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Locate the "multi fetch" heading paragraph and the paragraph right
#    after it ("Performs a git fetch on all dependencies.") and remove both,
#    together with the paragraph mark that ends the second one.
# ---------------------------------------------------------------------------
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "multi fetch`r") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find the 'multi fetch' heading paragraph."
}

$descIndex = $headingIndex + 1
$descPara = $d.Paragraphs.Item($descIndex)
if ($descPara.Range.Text -ne "Performs a git fetch on all dependencies.`r") {
    throw "Unexpected paragraph after 'multi fetch' heading."
}

$startPos = $d.Paragraphs.Item($headingIndex).Range.Start
$endPos = $descPara.Range.End
$removeRange = $d.Range($startPos, $endPos)
$removeRange.Delete()

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark. In the original document it sits right
#    after the "]" that ends the "multi branch ..." heading paragraph; it
#    needs to move to the end of the "... then installs it." paragraph
#    (right before the "multi install [--hooks]" heading).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*then installs it.`r") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'then installs it.' paragraph."
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$insertPos = $targetPara.Range.End - 1

# Placing a bookmark exactly at "paragraph end minus one" (i.e. right before
# the paragraph mark) in one step is unreliable, so nudge the position by
# temporarily inserting a placeholder character, anchoring the bookmark
# right before it, and then removing the placeholder again. Bookmarks track
# the surrounding text, so the bookmark stays put (now right before the
# paragraph mark) once the placeholder is gone.
$placeholderRange = $d.Range($insertPos, $insertPos)
$placeholderRange.InsertAfter("X")

$bookmarkRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$deleteRange = $d.Range($insertPos, $insertPos + 1)
$deleteRange.Delete()
